$d = $word.ActiveDocument

# 1. "theses ambition" -> "these ambitions"  (fix spelling + pluralize)
$d.Content.Find.Execute("theses ambition,", $true, $false, $false, $false, $false, $true, 1, $false, "these ambitions,", 2)

# 2. "We have done the best as we could.)" -> "We have done the best as we could for this analysis. )"
$d.Content.Find.Execute("We have done the best as we could.)", $true, $false, $false, $false, $false, $true, 1, $false, "We have done the best as we could for this analysis. )", 2)

# 3. "not straight. The median" -> "not straight). The median"
$d.Content.Find.Execute("not straight. The median", $true, $false, $false, $false, $false, $true, 1, $false, "not straight). The median", 2)

# 4. "are so different. " -> "are so different from state to state. "
$d.Content.Find.Execute("are so different. ", $true, $false, $false, $false, $false, $true, 1, $false, "are so different from state to state. ", 2)

# The two replacements above (1 and 4) each span across a proofErr spell/grammar
# marker that Word's proofing pass would normally retire once the flagged text
# is corrected/merged. Find/Replace leaves the now-orphaned marker behind, so
# clean those two stray tags up via a raw OOXML patch (same technique as
# ActiveDocument.WordOpenXML round-tripping).
$xml = $d.WordOpenXML
$xml = $xml.Replace('<w:proofErr w:type="spellStart"/><w:r><w:t xml:space="preserve">these ambitions, I can spend hours to hours </w:t></w:r>', '<w:r><w:t xml:space="preserve">these ambitions, I can spend hours to hours </w:t></w:r>')
$xml = $xml.Replace('<w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:color w:val="1D1C1D"/></w:rPr><w:t xml:space="preserve">are so different from state to state. </w:t></w:r>', '<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cstheme="minorHAnsi"/><w:color w:val="1D1C1D"/></w:rPr><w:t xml:space="preserve">are so different from state to state. </w:t></w:r>')
$d.WordOpenXML = $xml

Write-Host "Done"
